# Generate Report for Archive
# - Flip the localization status from "Ready for handoff" to "In Translation"
#   (this string is shared across the Overview summary row and each
#   language sheet's detail row).
# - The Status column narrows accordingly (report regenerated with the
#   shorter status text), on the Overview sheet (cols E/F) and on each
#   language sheet (col C).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update status text everywhere it appears -----------------------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value     = "In Translation"
$dede.Range("C2").Value     = "In Translation"

# --- Narrow the (now shorter) status columns -------------------------------
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth     = 12.5
$dede.Columns.Item(3).ColumnWidth     = 12.5
